$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for rows 2-24 that needs to
# move forward by one day (46074 -> 46075, i.e. 2026-02-21 -> 2026-02-22).
for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
